$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '44.441.09'
$ws.Range('E2').Value = '  +1.33%  '

# Row 3
$ws.Range('D3').Value = '2.249.33'
$ws.Range('E3').Value = '  +0.92%  '

# Row 4
$ws.Range('E4').Value = '  +0.22%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range('D5').Value = '308.78'
$ws.Range("D5").ClearFormats()
$ws.Range('E5').Value = '  +1.98%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range('D6').Value = '95.22'
$ws.Range("D6").ClearFormats()
$ws.Range('E6').Value = '  +2.20%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range('D7').Value = '0.573'
$ws.Range("D7").ClearFormats()
$ws.Range('E7').Value = '  +1.49%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range('D8').Value = '1.01'
$ws.Range("D8").ClearFormats()
$ws.Range('E8').Value = '  +0.17%  '

# Row 9
$ws.Range('E9').Value = '  +2.47%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range('D10').Value = '35.69'
$ws.Range("D10").ClearFormats()
$ws.Range('E10').Value = '  +5.56%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range('D11').Value = '0.0814'
$ws.Range("D11").ClearFormats()
$ws.Range('E11').Value = '  +2.66%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range('D12').Value = '7.26'
$ws.Range("D12").ClearFormats()
$ws.Range('E12').Value = '  +3.21%  '

# Row 13
$ws.Range('E13').Value = '  +1.92%  '

# Row 14
$ws.Range('D14').Value = '2.399.65'
$ws.Range('E14').Value = '  +6.05%  '

# Row 15
$ws.Range('E15').Value = '  +4.55%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range('D16').Value = '13.72'
$ws.Range("D16").ClearFormats()
$ws.Range('E16').Value = '  +2.70%  '

# Row 17
$ws.Range('D17').Value = '44.195.60'
$ws.Range('E17').Value = '  +1.20%  '

# Row 18
$ws.Range('D18').Value = '0.0₃0969'
$ws.Range('E18').Value = '  +2.29%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range('D19').Value = '12.34'
$ws.Range("D19").ClearFormats()
$ws.Range('E19').Value = '  +2.04%  '

# Row 20
$ws.Range('E20').Value = '  +5.40%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range('D21').Value = '66.20'
$ws.Range("D21").ClearFormats()
$ws.Range('E21').Value = '  +3.43%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range('D22').Value = '3.06'
$ws.Range("D22").ClearFormats()
$ws.Range('E22').Value = '  +6.07%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range('D23').Value = '238.31'
$ws.Range("D23").ClearFormats()
$ws.Range('E23').Value = '  +1.94%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range('D24').Value = '2.02'
$ws.Range("D24").ClearFormats()
$ws.Range('E24').Value = '  +6.44%  '

# Row 25
$ws.Range('E25').Value = '  +0.13%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range('D26').Value = '2.24'
$ws.Range("D26").ClearFormats()
$ws.Range('E26').Value = '  +6.35%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range('D27').Value = '38.42'
$ws.Range("D27").ClearFormats()
$ws.Range('E27').Value = '  +7.48%  '

# Row 28
$ws.Range('E28').Value = '  +1.93%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range('D29').Value = '6.02'
$ws.Range("D29").ClearFormats()
$ws.Range('E29').Value = '  +2.84%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range('D30').Value = '20.16'
$ws.Range("D30").ClearFormats()
$ws.Range('E30').Value = '  +2.53%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range('D31').Value = '154.60'
$ws.Range("D31").ClearFormats()
$ws.Range('E31').Value = '  +2.25%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range('D32').Value = '0.0802'
$ws.Range("D32").ClearFormats()
$ws.Range('E32').Value = '  +1.27%  '

# Row 33
$ws.Range('E33').Value = '  +1.01%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range('D34').Value = '3.14'
$ws.Range("D34").ClearFormats()
$ws.Range('E34').Value = '  -2.40%  '

# Row 35
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").NumberFormat = "@"
$ws.Range('D35').Value = '0.120'
$ws.Range("D35").ClearFormats()
$ws.Range('E35').Value = '  +2.77%  '

# Row 36
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").NumberFormat = "@"
$ws.Range('D36').Value = '0.109'
$ws.Range("D36").ClearFormats()
$ws.Range('E36').Value = '  +4.36%  '

# Row 37
$ws.Range('E37').Value = '  +6.18%  '

# Row 38
$ws.Range('E38').Value = '  +8.46%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range('D39').Value = '14.62'
$ws.Range("D39").ClearFormats()
$ws.Range('E39').Value = '  +2.64%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range('D40').Value = '3.85'
$ws.Range("D40").ClearFormats()
$ws.Range('E40').Value = '  +3.35%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range('D41').Value = '0.0305'
$ws.Range("D41").ClearFormats()
$ws.Range('E41').Value = '  +3.51%  '

# Row 42
$ws.Range('E42').Value = '  +0.23%  '

# Row 43
$ws.Range('D43').Value = '1.750.34'
$ws.Range('E43').Value = '  +1.30%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range('D44').Value = '0.195'
$ws.Range("D44").ClearFormats()
$ws.Range('E44').Value = '  +6.73%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range('D45').Value = '81.25'
$ws.Range("D45").ClearFormats()
$ws.Range('E45').Value = '  -1.56%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range('D46').Value = '71.51'
$ws.Range("D46").ClearFormats()
$ws.Range('E46').Value = '  +6.52%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range('D47').Value = '100.20'
$ws.Range("D47").ClearFormats()
$ws.Range('E47').Value = '  +1.53%  '

# Row 48
$ws.Range('B48').Value = 'THORChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D48").NumberFormat = "@"
$ws.Range('D48').Value = '4.95'
$ws.Range("D48").ClearFormats()
$ws.Range('E48').Value = '  +1.47%  '

# Row 49
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range('D49').Value = '1.61'
$ws.Range("D49").ClearFormats()
$ws.Range('E49').Value = '  +8.81%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range('D50').Value = '56.03'
$ws.Range("D50").ClearFormats()
$ws.Range('E50').Value = '  +5.13%  '

# Row 51
$ws.Range('B51').Value = 'FraxShare'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D51").NumberFormat = "@"
$ws.Range('D51').Value = '8.20'
$ws.Range("D51").ClearFormats()
$ws.Range('E51').Value = '  +2.13%  '
